$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column U (04-jul) with its 24 hourly
# values, mirroring the header style of the last existing day column (T).
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting from T1 onto U1, then set its text.
$wsSpot.Range("T1").Copy()
$wsSpot.Range("U1").PasteSpecial(-4122)
$wsSpot.Range("U1").Value = "04-jul"

$wsSpot.Range("U2").Value = 97.8
$wsSpot.Range("U3").Value = 83
$wsSpot.Range("U4").Value = 85.37
$wsSpot.Range("U5").Value = 75.8
$wsSpot.Range("U6").Value = 67.39
$wsSpot.Range("U7").Value = 78.65000000000001
$wsSpot.Range("U8").Value = 83.59
$wsSpot.Range("U9").Value = 103.48
$wsSpot.Range("U10").Value = 100
$wsSpot.Range("U11").Value = 91.15000000000001
$wsSpot.Range("U12").Value = 70.01000000000001
$wsSpot.Range("U13").Value = 20.64
$wsSpot.Range("U14").Value = 22.64
$wsSpot.Range("U15").Value = 12.34
$wsSpot.Range("U16").Value = 12.2
$wsSpot.Range("U17").Value = 25.2
$wsSpot.Range("U18").Value = 21.88
$wsSpot.Range("U19").Value = 62.04
$wsSpot.Range("U20").Value = 90.02
$wsSpot.Range("U21").Value = 108.5
$wsSpot.Range("U22").Value = 111.6
$wsSpot.Range("U23").Value = 109.13
$wsSpot.Range("U24").Value = 111.8
$wsSpot.Range("U25").Value = 96.25

# ---------------------------------------------------------------------------
# Sheet "Gaz": append a new row 18 for 2025-07-02.
# The date column stores plain text (not real dates) in this workbook, so
# the cell is forced to Text format before assignment (then the format
# override is dropped again) to stop Excel auto-converting the ISO-looking
# string into a date serial number.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A18").NumberFormat = "@"
$wsGaz.Range("A18").Value = "2025-07-02"
$wsGaz.Range("A18").Style = "Normal"
$wsGaz.Range("B18").Value = 32.675

# ---------------------------------------------------------------------------
# Sheet "CO2": append a new row 18 for 2025-07-02.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A18").NumberFormat = "@"
$wsCo2.Range("A18").Value = "2025-07-02"
$wsCo2.Range("A18").Style = "Normal"
$wsCo2.Range("B18").Value = 71.40000000000001
